$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its original text formatting so numeric-looking
# price strings (e.g. "1.000", "0.9995") are stored as text, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.520.88"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.836.53"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "318.98"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "0.5317"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("D8").Value = "0.4030"
$ws.Range("E8").Value = "  +6.71%  "
$ws.Range("D9").Value = "0.07582"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").Value = "41.94"
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("D11").Value = "1.112"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "6.326"
$ws.Range("E12").Value = "  +2.77%  "
$ws.Range("D13").Value = "7.617"
$ws.Range("E13").Value = "  +4.38%  "
$ws.Range("D14").Value = "0.9990"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "20.86"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").Value = "1.834.47"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").Value = "89.89"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "0.00001073"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").Value = "0.06592"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").Value = "17.67"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").Value = "0.9985"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("D23").Value = "28.555.74"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "11.31"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").Value = "2.102"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").Value = "2.461"
$ws.Range("E26").Value = "  +6.44%  "
$ws.Range("D27").Value = "20.61"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").Value = "157.11"
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").Value = "2.046.78"
$ws.Range("E29").Value = "  +2.27%  "
$ws.Range("D30").Value = "124.11"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").Value = "1.131"
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("E32").Value = "  +4.78%  "
$ws.Range("D33").Value = "5.708"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("D34").Value = "3.653"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "0.07196"
$ws.Range("E35").Value = "  +11.55%  "
$ws.Range("D36").Value = "0.2269"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D39").Value = "8.808"
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("D40").Value = "11.39"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("D41").Value = "0.6291"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("D42").Value = "1.201"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("D44").Value = "0.9983"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("D46").Value = "3.711"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("D47").Value = "0.5855"
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("D48").Value = "125.84"
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").Value = "1.997"
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("D50").Value = "1.197"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "0.06927"
$ws.Range("E51").Value = "  +0.76%  "

# Rows 37/38: VeChain and InternetComputer(DFINITY) swap list positions
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02349"
$ws.Range("E37").Value = "  +2.47%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "5.258"
$ws.Range("E38").Value = "  +4.76%  "
